$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column letters whose values change for each data row (B..N minus E,J,K,L which stay 0)
$cols = @("B", "C", "D", "F", "G", "H", "I", "M", "N")

$data = @{
    2 = @(0.5909449888591212, 0.1858151788494524, 0.007537871369759586, 0.5825610311031255, 0.4371995582342123, 0.5295738793486322, 0.3845158688133452, 0.8948114962484226, 0.9364384619600088)
    3 = @(0.51659482688612, 0.1635262707360425, 0.006983210171679133, 0.5654755118877688, 0.4200445841948124, 0.5262639454938949, 0.3834650222636355, 0.7961706769635413, 0.9520463109238859)
    4 = @(0.47088022357309, 0.1497728220114141, 0.006640299512572057, 0.5554521569532014, 0.4099057404223885, 0.5245825432965887, 0.3831295463662308, 0.7361304351839948, 0.9620990368732816)
    5 = @(0.4522360348598795, 0.1441512032803018, 0.006499985187122803, 0.5514845221674065, 0.4058725713139779, 0.5239853939000625, 0.3830704012086237, 0.7117895922557977, 0.9663134124500177)
    6 = @(0.4491392922036539, 0.143216715827009, 0.006476651776218034, 0.55083274851075, 0.4052087955579111, 0.5238915489160973, 0.3830652554024709, 0.7077552495316013, 0.9670203147749374)
    7 = @(0.4706288416657856, 0.1496970754812423, 0.006638409495362652, 0.5553981751561068, 0.4098509498009264, 0.5245741337510736, 0.3831284350934077, 0.7358016637125644, 0.96215539674232)
    8 = @(0.5653226491328383, 0.1781441325703383, 0.007347119737339369, 0.5765726496017578, 0.4312022548507031, 0.528359657053656, 0.3840890513946036, 0.8606872994964334, 0.9417224731295706)
    9 = @(0.7504884749969847, 0.2333869869027012, 0.008717706346430987, 0.6218295593177459, 0.4762346415377579, 0.5385777866914907, 0.3884459153649829, 1.110050655036758, 0.9053884764408702)
    10 = @(0.8861858402748908, 0.2736442262276739, 0.009712302459526256, 0.6573972753243993, 0.5112979798911965, 0.5478049760166925, 0.3931769416572308, 1.296407209555724, 0.880983851241484)
    11 = @(0.947839609413279, 0.291887055360462, 0.01016194806390303, 0.6740903465190655, 0.5276897735240311, 0.5523797616332899, 0.3956663581971682, 1.381965682450755, 0.8703816917035585)
    12 = @(0.9711747867278291, 0.2987849564443081, 0.01033180120282751, 0.6804860241310564, 0.5339611955990051, 0.5541666329318673, 0.3966579149732965, 1.414484900020369, 0.8664390598138541)
    13 = @(0.966149676256407, 0.2972998286817585, 0.01029523911524421, 0.6791052833691964, 0.5326076651532219, 0.553779370422177, 0.3964421860638225, 1.407475844966683, 0.8672849607517934)
    14 = @(0.9497596494126128, 0.292454756311173, 0.01017593044502618, 0.6746150288381614, 0.5282044361267424, 0.5525256749499476, 0.3957469526227015, 1.384638608661703, 0.8700558805767979)
    15 = @(0.9397187259386328, 0.2894856689567007, 0.01010279557286253, 0.6718743213825462, 0.5255157148933449, 0.5517648550104184, 0.3953274773286424, 1.370666001517961, 0.8717625583271271)
    16 = @(0.8821550232393065, 0.2724505860049362, 0.009682859401859645, 0.6563167194590989, 0.5102356815270923, 0.5475136159863609, 0.3930210705485564, 1.290832161091984, 0.8816868091128027)
    17 = @(0.8468215485005999, 0.2619819857791867, 0.009424513935297085, 0.6469044768821561, 0.5009754346023385, 0.5450024255298871, 0.3916928066438032, 1.242062642584756, 0.8879032143755659)
    18 = @(0.8264916015270387, 0.2559541066606243, 0.009275657690899664, 0.6415390677567814, 0.4956906783613277, 0.5435935527959828, 0.3909605454099818, 1.214085042677539, 0.8915257333154791)
    19 = @(0.819607044207487, 0.2539120381368605, 0.009225212857000287, 0.6397307040098354, 0.4939084544627832, 0.5431226220055692, 0.3907180511363393, 1.204624682944939, 0.892760319547139)
    20 = @(0.8505835950462597, 0.263097072009856, 0.009452042556546303, 0.6479014268208374, 0.5019569039491074, 0.5452660704068251, 0.3918309169066205, 1.247246598316224, 0.8872366000072107)
    21 = @(0.9545741240294774, 0.2938781511487889, 0.01021098575439083, 0.675931902272751, 0.5294960227371348, 0.5528924346378119, 0.3959498302687905, 1.391343140058154, 0.8692400327548135)
    22 = @(1.022468947703885, 0.3139354836115729, 0.01070455630564737, 0.6946852202688802, 0.5478690175251302, 0.5581944836504817, 0.3989267976860091, 1.486222032108998, 0.8578990903072334)
    23 = @(0.9862388093249024, 0.3032360333663462, 0.01044135686928271, 0.6846363361216845, 0.5380284817947256, 0.5553355217757883, 0.3973117309848107, 1.435516535055058, 0.863913345431091)
    24 = @(0.8488828243860098, 0.2625929706105978, 0.009439597893635465, 0.6474505630957879, 0.501513059978592, 0.5451467680449724, 0.3917683795577034, 1.244902743585072, 0.8875378249328403)
    25 = @(0.7004550518394694, 0.2185002694737932, 0.008349051231949289, 0.6091823852605671, 0.4637084534734726, 0.5355126522094054, 0.3870002022819818, 1.042069775666079, 0.9148166156228115)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}
